# Implemented usecase for bank withdraw.
#
# Record the actual time spent ("Realno utroseno vreme", column C) for the
# three tasks belonging to usecase 3 ("Kao korisnik potrebno je da mogu da
# prebacim novac sa naloga/novcanika na povezani bankovni racun" - i.e. the
# bank withdraw usecase), rows 30-32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C30").Value = "5min"
$ws.Range("C31").Value = "5min"
$ws.Range("C32").Value = "5min"

# Leave the view scrolled/selected where the user was working.
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("C32").Select()
